$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching original inlineStr formatting),
# then restore the default "Normal" style so no stray number-format style sticks.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0670"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0984"
$ws.Range("D11").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.632"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0515"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.643"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0186"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0513"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.62"
$ws.Range("D50").Style = "Normal"

# Remaining cells are safe to set directly (non-ambiguous text / swapped text values).
$ws.Range("D2").Value = "34.488.81"
$ws.Range("D3").Value = "1.806.61"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  +4.67%  "
$ws.Range("E8").Value = "  +7.06%  "
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("D12").Value = "2.067.53"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  -3.87%  "
$ws.Range("D14").Value = "1.814.13"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "34.461.68"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "0.0₃0769"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("E23").Value = "  -1.97%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "1.310.01"
$ws.Range("E36").Value = "  -5.92%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +4.64%  "
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").Value = "1.968.20"
$ws.Range("E48").Value = "  -4.23%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -6.45%  "
